$wb = $excel.ActiveWorkbook

# Sheet "展览" - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5124
$ws1.Range("F4").Value = 9
$ws1.Range("F5").Value = 7423
$ws1.Range("F13").Value = 1751
$ws1.Range("F14").Value = 102
$ws1.Range("F16").Value = 2914
$ws1.Range("F20").Value = 498
$ws1.Range("F21").Value = 433
$ws1.Range("F22").Value = 454
$ws1.Range("F24").Value = 97
$ws1.Range("F25").Value = 1689
$ws1.Range("F26").Value = 1183
$ws1.Range("F27").Value = 91
$ws1.Range("F28").Value = 1377
$ws1.Range("F30").Value = 578
$ws1.Range("F31").Value = 23
$ws1.Range("F36").Value = 62
$ws1.Range("F37").Value = 2870
$ws1.Range("F40").Value = 62

# Sheet "全部类型" - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5124
$ws4.Range("F4").Value = 9
$ws4.Range("F5").Value = 7423
$ws4.Range("F12").Value = 4309
$ws4.Range("F13").Value = 1751
$ws4.Range("F14").Value = 102
$ws4.Range("F16").Value = 2914
$ws4.Range("F20").Value = 498
$ws4.Range("F21").Value = 433
$ws4.Range("F22").Value = 454
$ws4.Range("F25").Value = 97
$ws4.Range("F26").Value = 1689
$ws4.Range("F27").Value = 1183
$ws4.Range("F28").Value = 91
$ws4.Range("F29").Value = 1377
$ws4.Range("F31").Value = 578
$ws4.Range("F32").Value = 23
$ws4.Range("F37").Value = 62
$ws4.Range("F38").Value = 2870
$ws4.Range("F42").Value = 62
